$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "27.072.80"
$ws.Cells.Item(3, 4).Value = "1.675.43"
$ws.Cells.Item(3, 5).Value = "  +0.26%  "
$ws.Cells.Item(4, 5).Value = "  +0.06%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "215.34"
$ws.Cells.Item(5, 5).Value = "  +0.29%  "
$ws.Cells.Item(6, 5).Value = "  -0.09%  "
$ws.Cells.Item(7, 5).Value = "  +0.02%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.256"
$ws.Cells.Item(8, 5).Value = "  +1.84%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "21.24"
$ws.Cells.Item(9, 5).Value = "  +4.94%  "
$ws.Cells.Item(10, 5).Value = "  +0.28%  "
$ws.Cells.Item(11, 5).Value = "  -0.82%  "
$ws.Cells.Item(13, 4).Value = "1.677.97"
$ws.Cells.Item(13, 5).Value = "  +0.65%  "
$ws.Cells.Item(14, 5).Value = "  +0.99%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.535"
$ws.Cells.Item(15, 5).Value = "  +1.69%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "66.02"
$ws.Cells.Item(16, 5).Value = "  +0.82%  "
$ws.Cells.Item(17, 4).Value = "27.048.13"
$ws.Cells.Item(17, 5).Value = "  +0.47%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "8.17"
$ws.Cells.Item(18, 5).Value = "  +1.71%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "237.52"
$ws.Cells.Item(19, 5).Value = "  +1.99%  "
$ws.Cells.Item(20, 4).Value = "0.0₃0742"
$ws.Cells.Item(20, 5).Value = "  +1.37%  "
$ws.Cells.Item(21, 5).Value = "  +0.00%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "4.47"
$ws.Cells.Item(22, 5).Value = "  +1.00%  "
$ws.Cells.Item(23, 5).Value = "  +2.28%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "2.14"
$ws.Cells.Item(24, 5).Value = "  -2.01%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "146.36"
$ws.Cells.Item(25, 5).Value = "  +0.48%  "
$ws.Cells.Item(26, 5).Value = "  +1.62%  "
$ws.Cells.Item(27, 5).Value = "  +2.84%  "
$ws.Cells.Item(28, 5).Value = "  +0.39%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "0.999"
$ws.Cells.Item(29, 5).Value = "  -0.17%  "
$ws.Cells.Item(30, 5).Value = "  +0.09%  "
$ws.Cells.Item(31, 5).Value = "  +0.11%  "
$ws.Cells.Item(32, 5).Value = "  +0.87%  "
$ws.Cells.Item(33, 4).Value = "1.549.13"
$ws.Cells.Item(33, 5).Value = "  +6.47%  "
$ws.Cells.Item(34, 5).Value = "  +2.28%  "
$ws.Cells.Item(35, 5).Value = "  +3.80%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.598"
$ws.Cells.Item(36, 5).Value = "  +3.46%  "
$ws.Cells.Item(37, 5).Value = "  -1.05%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.924"
$ws.Cells.Item(38, 5).Value = "  +2.70%  "
$ws.Cells.Item(39, 5).Value = "  +2.17%  "
$ws.Cells.Item(40, 5).Value = "  +1.70%  "
$ws.Cells.Item(41, 5).Value = "  +0.04%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "67.61"
$ws.Cells.Item(42, 5).Value = "  +1.99%  "
$ws.Cells.Item(43, 5).Value = "  -2.24%  "
$ws.Cells.Item(44, 5).Value = "  -1.72%  "
$ws.Cells.Item(45, 4).Value = "1.819.73"
$ws.Cells.Item(45, 5).Value = "  +0.67%  "
$ws.Cells.Item(46, 5).Value = "  +0.58%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "90.93"
$ws.Cells.Item(47, 5).Value = "  +0.59%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "1.57"
$ws.Cells.Item(48, 5).Value = "  +2.35%  "
$ws.Cells.Item(49, 5).Value = "  +1.72%  "
$ws.Cells.Item(50, 5).Value = "  +2.52%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "8.03"
$ws.Cells.Item(51, 5).Value = "  +5.46%  "
